# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Julio de 2020 a las 19:33"

# --- Re-rank swapped country labels (values identical, only the display order/name changed) ---
# Austria / Moldavia swap position (row 60/61) - also gets new daily figures below
$ws.Range("A60").Value = "Moldavia"
$ws.Range("A61").Value = "Austria"

# Seychelles / Lesoto swap position (row 184/185) - figures unchanged
$ws.Range("A184").Value = "Lesoto"
$ws.Range("A185").Value = "Seychelles"

# Islas Malvinas / Groenlandia swap position (row 209/210) - figures unchanged
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- Update numeric columns (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 3184938
$ws.Range("C4").Value = 26006
$ws.Range("D4").Value = 1403684
$ws.Range("E4").Value = 1645973
$ws.Range("G4").Value = 419
$ws.Range("H4").Value = 135281

# Row 5 - Brasil
$ws.Range("B5").Value = 1727279
$ws.Range("C5").Value = 11083
$ws.Range("E5").Value = 506457
$ws.Range("G5").Value = 300
$ws.Range("H5").Value = 68355

# Row 6 - India
$ws.Range("B6").Value = 794196
$ws.Range("C6").Value = 25144
$ws.Range("D6").Value = 495895
$ws.Range("E6").Value = 276679
$ws.Range("G6").Value = 478
$ws.Range("H6").Value = 21622

# Row 10 - España
$ws.Range("B10").Value = 300136
$ws.Range("C10").Value = 543
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 28401

# Row 18 - Turquia
$ws.Range("B18").Value = 209962
$ws.Range("C18").Value = 1024
$ws.Range("D18").Value = 190390
$ws.Range("E18").Value = 14272
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 5300

# Row 19 - Alemania
$ws.Range("B19").Value = 198926
$ws.Range("C19").Value = 161
$ws.Range("E19").Value = 6203
$ws.Range("G19").Value = 8
$ws.Range("H19").Value = 9123

# Row 47 - Israel
$ws.Range("B47").Value = 34825
$ws.Range("C47").Value = 1268
$ws.Range("D47").Value = 18452
$ws.Range("E47").Value = 16025
$ws.Range("G47").Value = 4
$ws.Range("H47").Value = 348

# Row 60 - now Moldavia (new figures)
$ws.Range("B60").Value = 18666
$ws.Range("C60").Value = 195
$ws.Range("D60").Value = 11936
$ws.Range("E60").Value = 6106
$ws.Range("G60").Value = 10
$ws.Range("H60").Value = 624

# Row 61 - now Austria (former Austria figures)
$ws.Range("B61").Value = 18615
$ws.Range("C61").Value = 102
$ws.Range("D61").Value = 16758
$ws.Range("E61").Value = 1151
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 706

# Row 65 - Marruecos
$ws.Range("B65").Value = 15079
$ws.Range("C65").Value = 308
$ws.Range("D65").Value = 11447
$ws.Range("E65").Value = 3390

# Row 93 - Estado de Palestina
$ws.Range("E93").Value = 4671
$ws.Range("G93").Value = 4
$ws.Range("H93").Value = 24

# Row 96 - Luxemburgo
$ws.Range("B96").Value = 4719
$ws.Range("C96").Value = 69
$ws.Range("E96").Value = 553

# Row 111 - Mali
$ws.Range("B111").Value = 2370
$ws.Range("C111").Value = 12
$ws.Range("D111").Value = 1621
$ws.Range("E111").Value = 629

# Row 112 - Sri Lanka
$ws.Range("B112").Value = 2154
$ws.Range("C112").Value = 60
$ws.Range("E112").Value = 164

# Row 125 - Cabo Verde
$ws.Range("B125").Value = 1552
$ws.Range("C125").Value = 10
$ws.Range("E125").Value = 804

# Row 130 - Tunez
$ws.Range("B130").Value = 1231
$ws.Range("C130").Value = 10
$ws.Range("D130").Value = 1055
$ws.Range("E130").Value = 126

# Row 135 - Suazilandia
$ws.Range("E135").Value = 534
$ws.Range("G135").Value = 2
$ws.Range("H135").Value = 16

# Row 137 - Mozambique
$ws.Range("B137").Value = 1092
$ws.Range("C137").Value = 21
$ws.Range("D137").Value = 340
$ws.Range("E137").Value = 743
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 9
